# Update the "想去人数" (wish-to-attend count) numbers in the F column
# across the workbook's sheets, matching the regenerated data snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 2453
$ws1.Range("F3").Value = 720
$ws1.Range("F6").Value = 685
$ws1.Range("F8").Value = 887
$ws1.Range("F9").Value = 564
$ws1.Range("F13").Value = 437
$ws1.Range("F16").Value = 1067
$ws1.Range("F17").Value = 23900
$ws1.Range("F18").Value = 2199
$ws1.Range("F21").Value = 29
$ws1.Range("F24").Value = 207
$ws1.Range("F29").Value = 44
$ws1.Range("F30").Value = 340
$ws1.Range("F32").Value = 432

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value = 19
$ws2.Range("F9").Value = 247
$ws2.Range("F19").Value = 4113

# --- Sheet "全部类型" (All types, aggregated view) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 2453
$ws4.Range("F6").Value = 720
$ws4.Range("F9").Value = 685
$ws4.Range("F16").Value = 887
$ws4.Range("F17").Value = 564
$ws4.Range("F20").Value = 437
$ws4.Range("F23").Value = 1067
$ws4.Range("F24").Value = 23900
$ws4.Range("F25").Value = 19
$ws4.Range("F26").Value = 247
$ws4.Range("F30").Value = 2199
$ws4.Range("F33").Value = 29
$ws4.Range("F37").Value = 207
$ws4.Range("F43").Value = 44
$ws4.Range("F47").Value = 432
